$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Account Name" / "Phone" request fields below the existing
# ContactName / Email / Subject rows, then update the contact's details to
# reflect the new request (Rose Mwangi's account creation), in the same
# order the data was actually typed into the sheet.
$ws.Range("A4").Value = "Account Name"
$ws.Range("A5").Value = "Phone"
$ws.Range("B3").Value = "Account creation"
$ws.Range("B1").Value = "Rose"
$ws.Range("B4").Value = "Rose Mwangi"
$ws.Range("B2").Value = "rm.4@gmail.com"
$ws.Range("B5").Value = 727290683

# Widen column A to fit the new, longer labels
$ws.Columns.Item(1).ColumnWidth = 13.1

# Move the active selection to the row below the newly entered data
$ws.Range("A6:XFD6").Select()
